$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1821086261980831
$ws.Range("C2").Value = 0.5718849840255591
$ws.Range("J2").Value = 0.01597444089456869
$ws.Range("P2").Value = 0.1469648562300319
$ws.Range("S2").Value = 0.08306709265175719
$ws.Range("B3").Value = 0.00546448087431694
$ws.Range("C3").Value = 0.00546448087431694
$ws.Range("J3").Value = 0.01639344262295082
$ws.Range("P3").Value = 0.7978142076502732
$ws.Range("S3").Value = 0.1748633879781421
$ws.Range("J4").Value = 0.1272727272727273
$ws.Range("P4").Value = 0.5272727272727272
$ws.Range("S4").Value = 0.3454545454545455
$ws.Range("B6").Value = 0.07526881720430108
$ws.Range("D6").Value = 0.01612903225806452
$ws.Range("F6").Value = 0.03225806451612903
$ws.Range("J6").Value = 0.2311827956989247
$ws.Range("O6").Value = 0.02688172043010753
$ws.Range("Q6").Value = 0.1935483870967742
$ws.Range("R6").Value = 0.09139784946236559
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.1279620853080569
$ws.Range("D7").Value = 0.01895734597156398
$ws.Range("E7").Value = 0.004739336492890996
$ws.Range("F7").Value = 0.04739336492890995
$ws.Range("J7").Value = 0.1137440758293839
$ws.Range("O7").Value = 0.01421800947867299
$ws.Range("Q7").Value = 0.1706161137440758
$ws.Range("R7").Value = 0.04739336492890995
$ws.Range("S7").Value = 0.4549763033175355
$ws.Range("B8").Value = 0.1217183770883055
$ws.Range("D8").Value = 0.01909307875894988
$ws.Range("F8").Value = 0.03579952267303103
$ws.Range("J8").Value = 0.1193317422434367
$ws.Range("O8").Value = 0.02147971360381861
$ws.Range("Q8").Value = 0.1933174224343675
$ws.Range("R8").Value = 0.05966587112171837
$ws.Range("S8").Value = 0.4295942720763723
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.02857142857142857
$ws.Range("F9").Value = 0.05
$ws.Range("J9").Value = 0.1357142857142857
$ws.Range("O9").Value = 0.007142857142857143
$ws.Range("Q9").Value = 0.1642857142857143
$ws.Range("R9").Value = 0.05714285714285714
$ws.Range("S9").Value = 0.4571428571428571
$ws.Range("B10").Value = 0.1279373368146214
$ws.Range("D10").Value = 0.03220191470844212
$ws.Range("E10").Value = 0.0008703220191470844
$ws.Range("F10").Value = 0.0670147954743255
$ws.Range("J10").Value = 0.1044386422976501
$ws.Range("O10").Value = 0.01827676240208877
$ws.Range("Q10").Value = 0.2149695387293299
$ws.Range("R10").Value = 0.05395996518711924
$ws.Range("S10").Value = 0.3803307223672759
$ws.Range("G11").Value = 0.1585365853658537
$ws.Range("J11").Value = 0.1097560975609756
$ws.Range("K11").Value = 0.2195121951219512
$ws.Range("L11").Value = 0.5
$ws.Range("S11").Value = 0.01219512195121951
$ws.Range("G12").Value = 0.7784431137724551
$ws.Range("J12").Value = 0.1736526946107785
$ws.Range("L12").Value = 0.02395209580838323
$ws.Range("S12").Value = 0.02395209580838323
$ws.Range("G13").Value = 0.6379310344827587
$ws.Range("J13").Value = 0.3620689655172414
$ws.Range("F15").Value = 0.03092783505154639
$ws.Range("H15").Value = 0.1701030927835052
$ws.Range("I15").Value = 0.04639175257731959
$ws.Range("J15").Value = 0.3814432989690721
$ws.Range("K15").Value = 0.06185567010309279
$ws.Range("M15").Value = 0.0154639175257732
$ws.Range("O15").Value = 0.03608247422680412
$ws.Range("S15").Value = 0.2577319587628866
$ws.Range("F16").Value = 0.03669724770642202
$ws.Range("H16").Value = 0.1972477064220184
$ws.Range("I16").Value = 0.04587155963302753
$ws.Range("J16").Value = 0.4036697247706422
$ws.Range("K16").Value = 0.1009174311926606
$ws.Range("M16").Value = 0.02293577981651376
$ws.Range("O16").Value = 0.07339449541284404
$ws.Range("S16").Value = 0.1192660550458716
$ws.Range("F17").Value = 0.01428571428571429
$ws.Range("H17").Value = 0.1714285714285714
$ws.Range("I17").Value = 0.08333333333333333
$ws.Range("J17").Value = 0.4285714285714285
$ws.Range("K17").Value = 0.1119047619047619
$ws.Range("M17").Value = 0.02857142857142857
$ws.Range("O17").Value = 0.05476190476190476
$ws.Range("S17").Value = 0.1071428571428571
$ws.Range("F18").Value = 0.008264462809917356
$ws.Range("H18").Value = 0.2148760330578512
$ws.Range("I18").Value = 0.06611570247933884
$ws.Range("J18").Value = 0.3636363636363636
$ws.Range("K18").Value = 0.1322314049586777
$ws.Range("M18").Value = 0.04958677685950413
$ws.Range("O18").Value = 0.08264462809917356
$ws.Range("S18").Value = 0.08264462809917356
$ws.Range("F19").Value = 0.02451394759087067
$ws.Range("H19").Value = 0.2054099746407439
$ws.Range("I19").Value = 0.06931530008453085
$ws.Range("J19").Value = 0.3575655114116653
$ws.Range("K19").Value = 0.1318681318681319
$ws.Range("M19").Value = 0.02789518174133559
$ws.Range("O19").Value = 0.0676246830092984
$ws.Range("S19").Value = 0.1158072696534235
